# This script updates the NATMI LR-pairs TPM values in the Cd34-Selp sheet,
# replacing stale expression/specificity figures with the newly computed
# TPM-based values, per the commit "update scripts wuth new tpm".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 53.64296733333333
$ws.Range("H2").Value = 160.928902
$ws.Range("I2").Value = 0.1405570576660657
$ws.Range("J2").Value = 0.1405570576660657
$ws.Range("M2").Value = 1.343359
$ws.Range("N2").Value = 4.030077
$ws.Range("O2").Value = 0.736296379391111
$ws.Range("P2").Value = 0.7362963793911109
$ws.Range("Q2").Value = 72.06176295393934
$ws.Range("R2").Value = 648.5558665854541
$ws.Range("S2").Value = 0.1034916526573918
$ws.Range("T2").Value = 0.1034916526573918

# Row 3
$ws.Range("G3").Value = 53.64296733333333
$ws.Range("H3").Value = 160.928902
$ws.Range("I3").Value = 0.1405570576660657
$ws.Range("J3").Value = 0.1405570576660657
$ws.Range("O3").Value = 0.1764523396969075
$ws.Range("P3").Value = 0.1764523396969075
$ws.Range("Q3").Value = 17.26949504548933
$ws.Range("R3").Value = 155.425455409404
$ws.Range("S3").Value = 0.02480162168609044
$ws.Range("T3").Value = 0.02480162168609044

# Row 4
$ws.Range("G4").Value = 53.64296733333333
$ws.Range("H4").Value = 160.928902
$ws.Range("I4").Value = 0.1405570576660657
$ws.Range("J4").Value = 0.1405570576660657
$ws.Range("O4").Value = 0.08725128091198156
$ws.Range("P4").Value = 0.08725128091198156
$ws.Range("Q4").Value = 8.539334564847778
$ws.Range("R4").Value = 76.85401108363
$ws.Range("S4").Value = 0.01226378332258349
$ws.Range("T4").Value = 0.01226378332258349

# Row 5
$ws.Range("I5").Value = 0.83973167405618
$ws.Range("J5").Value = 0.8397316740561799
$ws.Range("M5").Value = 1.343359
$ws.Range("N5").Value = 4.030077
$ws.Range("O5").Value = 0.736296379391111
$ws.Range("P5").Value = 0.7362963793911109
$ws.Range("Q5").Value = 430.5194334994994
$ws.Range("R5").Value = 3874.674901495494
$ws.Range("S5").Value = 0.6182913912676019
$ws.Range("T5").Value = 0.6182913912676017

# Row 6
$ws.Range("I6").Value = 0.83973167405618
$ws.Range("J6").Value = 0.8397316740561799
$ws.Range("O6").Value = 0.1764523396969075
$ws.Range("P6").Value = 0.1764523396969075
$ws.Range("S6").Value = 0.1481726186048138
$ws.Range("T6").Value = 0.1481726186048138

# Row 7
$ws.Range("I7").Value = 0.83973167405618
$ws.Range("J7").Value = 0.8397316740561799
$ws.Range("O7").Value = 0.08725128091198156
$ws.Range("P7").Value = 0.08725128091198156
$ws.Range("S7").Value = 0.07326766418376429
$ws.Range("T7").Value = 0.07326766418376429

# Row 8
$ws.Range("G8").Value = 7.522716666666668
$ws.Range("I8").Value = 0.01971126827775425
$ws.Range("J8").Value = 0.01971126827775425
$ws.Range("M8").Value = 1.343359
$ws.Range("N8").Value = 4.030077
$ws.Range("O8").Value = 0.736296379391111
$ws.Range("P8").Value = 0.7362963793911109
$ws.Range("Q8").Value = 10.10570913861667
$ws.Range("R8").Value = 90.95138224755001
$ws.Range("S8").Value = 0.01451333546611732
$ws.Range("T8").Value = 0.01451333546611731

# Row 9
$ws.Range("G9").Value = 7.522716666666668
$ws.Range("I9").Value = 0.01971126827775425
$ws.Range("J9").Value = 0.01971126827775425
$ws.Range("O9").Value = 0.1764523396969075
$ws.Range("P9").Value = 0.1764523396969075
$ws.Range("S9").Value = 0.003478099406003169
$ws.Range("T9").Value = 0.003478099406003169

# Row 10
$ws.Range("G10").Value = 7.522716666666668
$ws.Range("I10").Value = 0.01971126827775425
$ws.Range("J10").Value = 0.01971126827775425
$ws.Range("O10").Value = 0.08725128091198156
$ws.Range("P10").Value = 0.08725128091198156
$ws.Range("S10").Value = 0.001719833405633767
$ws.Range("T10").Value = 0.001719833405633767
